# Audit_Template.xlsx - "Fixed all Auto Complete on Forms"
#
# 1) Clear the sample/demo data that used to live in rows 2-4 (Case ID,
#    SAP numbers, dates, Pass/Fail, Qulaity*, SLA/Non SLA samples, etc.)
#    but keep all the cell formatting (borders/fills) intact.
# 2) Add two new header cells describing the "Note" guidance for filling
#    the form: L1 = "Note:" and M1 = the instructional text.
# 3) Re-flow row 1 / columns to accommodate the new note columns.
# 4) Park the active selection on L5, matching the saved view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Wipe out the old sample rows, keep the styling ---------------------
$ws.Range("A2:K4").ClearContents()

# --- 2. New "Note" header cells --------------------------------------------
# Write the values in column order (L then M) so the shared-string table
# keeps the same ordering Excel produced originally.
$ws.Range("L1").Value = "Note:"
$ws.Range("M1").Value = "For BC, EUC, CC or Fatal. Please use ""Pass"" or ""Fail""."

# Format M1 first (long wrapped note) ...
$ws.Range("M1").WrapText = $true
$ws.Range("M1").HorizontalAlignment = -4131 ## xlLeft
$ws.Range("M1").VerticalAlignment = -4108   ## xlCenter

# ... then L1 ("Note:" label), so the new cell styles are appended to
# cellXfs in the same order Excel originally recorded them.
$ws.Range("L1").HorizontalAlignment = -4108 ## xlCenter
$ws.Range("L1").VerticalAlignment = -4160   ## xlTop

# --- 3. Row height / column widths -----------------------------------------
$ws.Rows.Item(1).RowHeight = 25.5
$ws.Columns.Item(12).ColumnWidth = 4.166666666666667   ## -> stored width 5
$ws.Columns.Item(13).ColumnWidth = 22.5                ## -> note column width

# --- 4. Restore the saved selection -----------------------------------------
$ws.Range("L5").Select()
